# feat: Add multiple functionalities
# 1. exception handling in all functions
# 2. logging in proper manner
# 3. selecting elements in proper manner

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab to a more descriptive name
$ws.Name = "Product Data"

# 2. Set explicit column widths for readability
# (ColumnWidth is in "characters"; Excel stores the saved <col width="..">
#  with a constant +5/6 character padding added on top of what is set here,
#  so we subtract that padding to land exactly on the target stored widths
#  of 35 / 15 / 168 / 21.)
$padding = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 35 - $padding
$ws.Columns.Item(2).ColumnWidth = 15 - $padding
$ws.Columns.Item(3).ColumnWidth = 168 - $padding
$ws.Columns.Item(4).ColumnWidth = 21 - $padding

# 3. Refresh the Timestamp column with the latest run time
$newTimestamp = "2025-06-26 14:16:26"
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 4).Value = $newTimestamp
}
